# Update "想去人数" (F column) figures across the three affected sheets
# (展览, 演出, 全部类型) to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1207
$ws1.Range("F3").Value  = 1105
$ws1.Range("F4").Value  = 871
$ws1.Range("F5").Value  = 95
$ws1.Range("F6").Value  = 57
$ws1.Range("F7").Value  = 630
$ws1.Range("F8").Value  = 77
$ws1.Range("F10").Value = 32
$ws1.Range("F11").Value = 2180
$ws1.Range("F12").Value = 1546
$ws1.Range("F13").Value = 1218
$ws1.Range("F16").Value = 477
$ws1.Range("F17").Value = 709
$ws1.Range("F18").Value = 258
$ws1.Range("F19").Value = 1078
$ws1.Range("F22").Value = 4083
$ws1.Range("F23").Value = 204
$ws1.Range("F24").Value = 140
$ws1.Range("F28").Value = 588
$ws1.Range("F29").Value = 19
$ws1.Range("F30").Value = 58
$ws1.Range("F31").Value = 32
$ws1.Range("F32").Value = 233
$ws1.Range("F33").Value = 351
$ws1.Range("F34").Value = 911
$ws1.Range("F36").Value = 85
$ws1.Range("F37").Value = 104
$ws1.Range("F38").Value = 102

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 771

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1207
$ws4.Range("F4").Value  = 771
$ws4.Range("F5").Value  = 1105
$ws4.Range("F6").Value  = 871
$ws4.Range("F9").Value  = 95
$ws4.Range("F10").Value = 57
$ws4.Range("F11").Value = 630
$ws4.Range("F12").Value = 77
$ws4.Range("F15").Value = 32
$ws4.Range("F16").Value = 2180
$ws4.Range("F17").Value = 1546
$ws4.Range("F18").Value = 1218
$ws4.Range("F21").Value = 477
$ws4.Range("F23").Value = 709
$ws4.Range("F24").Value = 258
$ws4.Range("F25").Value = 1078
$ws4.Range("F28").Value = 4083
$ws4.Range("F29").Value = 204
$ws4.Range("F30").Value = 140
$ws4.Range("F34").Value = 588
$ws4.Range("F35").Value = 19
$ws4.Range("F36").Value = 58
$ws4.Range("F37").Value = 32
$ws4.Range("F38").Value = 233
$ws4.Range("F39").Value = 351
$ws4.Range("F40").Value = 911
$ws4.Range("F42").Value = 85
$ws4.Range("F43").Value = 104
$ws4.Range("F44").Value = 102
